# Automating the Input files for Transavia
# Append the next day's hourly predictions (date 45374, hours 0-9) to the
# Production_Predictions sheet, extending the data range from A1:C87 to A1:C97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateSerial = 45374
$values = @(
    @(0, 0.007000000216066837),
    @(1, 0.007000000216066837),
    @(2, 0.007000000216066837),
    @(3, 0.007000000216066837),
    @(4, 0.007000000216066837),
    @(5, 0.007000000216066837),
    @(6, 0.06599999964237213),
    @(7, 0.8109999895095825),
    @(8, 1.927999973297119),
    @(9, 3.059000015258789)
)

$startRow = 88
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $hour = $values[$i][0]
    $pred = $values[$i][1]

    $ws.Cells.Item($row, 1).Value = $dateSerial
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 1, 1).Style
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

    $ws.Cells.Item($row, 2).Value = $hour

    $ws.Cells.Item($row, 3).Value = $pred
    $ws.Cells.Item($row, 3).Style = $ws.Cells.Item($row - 1, 3).Style
    $ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat
}
